# [github-234] Ensure the hours position is correct.
# Adds four new "days h" / "days h am/pm" TEXT() format test rows (43-46)
# to the "Tests" sheet, exercising the "d \d\a\y\s h" / d "days" h"
# (with and without am/pm) custom number-format patterns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests")

$serial = 17816.607951388887

# Row 43: d \d\a\y\s h
$ws.Cells.Item(42, 3).Copy()
$ws.Cells.Item(43, 3).PasteSpecial(-4122)
$ws.Cells.Item(43, 3).Value = $serial
$ws.Cells.Item(43, 2).Value = "d \d\a\y\s h"
$ws.Cells.Item(43, 2).NumberFormat = "@"
$ws.Cells.Item(43, 2).Font.Name = "Lucida Sans Regular"
$ws.Cells.Item(43, 1).Formula = "=TEXT(C43,B43)"
$ws.Cells.Item(43, 4).Value = "Time"

# Row 44: d "days" h
$ws.Cells.Item(42, 3).Copy()
$ws.Cells.Item(44, 3).PasteSpecial(-4122)
$ws.Cells.Item(44, 3).Value = $serial
$ws.Cells.Item(44, 2).Value = "d ""days"" h"
$ws.Cells.Item(44, 2).NumberFormat = "@"
$ws.Cells.Item(44, 2).Font.Name = "Lucida Sans Regular"
$ws.Cells.Item(44, 1).Formula = "=TEXT(C44,B44)"
$ws.Cells.Item(44, 4).Value = "Time"

# Row 45: d \d\a\y\s h a/p
$ws.Cells.Item(42, 3).Copy()
$ws.Cells.Item(45, 3).PasteSpecial(-4122)
$ws.Cells.Item(45, 3).Value = $serial
$ws.Cells.Item(45, 2).Value = "d \d\a\y\s h a/p"
$ws.Cells.Item(45, 2).NumberFormat = "@"
$ws.Cells.Item(45, 2).Font.Name = "Lucida Sans Regular"
$ws.Cells.Item(45, 1).Formula = "=TEXT(C45,B45)"
$ws.Cells.Item(45, 4).Value = "Time"

# Row 46: d "days" h am/pm
$ws.Cells.Item(42, 3).Copy()
$ws.Cells.Item(46, 3).PasteSpecial(-4122)
$ws.Cells.Item(46, 3).Value = $serial
$ws.Cells.Item(46, 2).Value = "d ""days"" h am/pm"
$ws.Cells.Item(46, 2).NumberFormat = "@"
$ws.Cells.Item(46, 2).Font.Name = "Lucida Sans Regular"
$ws.Cells.Item(46, 1).Formula = "=TEXT(C46,B46)"
$ws.Cells.Item(46, 4).Value = "Time"

$ws.Range("C46").Select()
